$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values formatted as plain text (with "." as
# thousands separators in some rows), so force text format before
# assigning to avoid Excel auto-converting them to numbers/dates.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.070.80"
$ws.Range("E2").Value = "  +3.73%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.419.96"
$ws.Range("E3").Value = "  +3.23%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "552.90"
$ws.Range("E5").Value = "  +2.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.69"
$ws.Range("E6").Value = "  +2.66%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +2.53%  "
$ws.Range("E9").Value = "  +3.73%  "
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("E12").Value = "  -1.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.87"
$ws.Range("E13").Value = "  +4.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.849.63"
$ws.Range("E14").Value = "  +3.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.974.77"
$ws.Range("E15").Value = "  +3.63%  "
$ws.Range("E16").Value = "  +3.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.422.40"
$ws.Range("E17").Value = "  +3.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.33"
$ws.Range("E18").Value = "  +6.36%  "
$ws.Range("E19").Value = "  +1.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "332.39"
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.30"
$ws.Range("E23").Value = "  +3.77%  "
$ws.Range("E24").Value = "  +3.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.59"
$ws.Range("E25").Value = "  +3.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.51%  "
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("E28").Value = "  +6.55%  "
$ws.Range("E29").Value = "  +1.28%  "
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.26"
$ws.Range("E30").Value = "  +2.25%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "169.19"
$ws.Range("E31").Value = "  -0.66%  "
$ws.Range("E32").Value = "  +2.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.67"
$ws.Range("E33").Value = "  +1.98%  "
$ws.Range("E35").Value = "  +5.73%  "
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("E37").Value = "  +0.95%  "
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.419"
$ws.Range("E39").Value = "  +11.05%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "39.44"
$ws.Range("E40").Value = "  +0.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "313.76"
$ws.Range("E41").Value = "  +8.72%  "
$ws.Range("E42").Value = "  +1.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "139.21"
$ws.Range("E43").Value = "  -1.42%  "
$ws.Range("E45").Value = "  +2.47%  "
$ws.Range("E46").Value = "  +2.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.415"
$ws.Range("E47").Value = "  +7.85%  "
$ws.Range("E48").Value = "  +1.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0226"
$ws.Range("E49").Value = "  +1.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.71"
$ws.Range("E50").Value = "  +2.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.04"
$ws.Range("E51").Value = "  -0.25%  "
